# "New Microsite scripts support to Beta server"
# Appends new sprint-run history rows to the AMSIN, BETA and AMS sheets
# (mirrors the existing row layout/format of each sheet), and backfills
# formatting on AMS!A30:G30 that the prior edit had left unstyled.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper-ish inline pattern used throughout:
#   1) stamp text columns with NumberFormat "@" first so values that
#      look like dates ("2022-09-16") are kept as literal text instead
#      of being auto-parsed into date serials;
#   2) write all the cell values;
#   3) copy the row immediately above and paste-special (formats only)
#      onto the new row so it inherits the sheet's established look
#      (date/time display on column B, etc.) instead of Excel defaults.
# ---------------------------------------------------------------------

# --- Sheet "AMSIN": three new rows (52-54) ----------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Range("A52").NumberFormat = "@"
$wsAmsin.Range("A52").Value = "2022-09-16"
$wsAmsin.Range("B52").Value = 44820.60906491898
$wsAmsin.Range("C52").NumberFormat = "@"
$wsAmsin.Range("C52").Value = "fstcyc167"
$wsAmsin.Range("D52").Value = 75
$wsAmsin.Range("E52").Value = 74
$wsAmsin.Range("F52").Value = 1
$wsAmsin.Range("G52").Value = 2.28

$wsAmsin.Range("A53").NumberFormat = "@"
$wsAmsin.Range("A53").Value = "2022-09-19"
$wsAmsin.Range("B53").Value = 44823.6235595949
$wsAmsin.Range("C53").NumberFormat = "@"
$wsAmsin.Range("C53").Value = "scndcycle167"
$wsAmsin.Range("D53").Value = 75
$wsAmsin.Range("E53").Value = 75
$wsAmsin.Range("F53").Value = 0
$wsAmsin.Range("G53").Value = 2.01

$wsAmsin.Range("A54").NumberFormat = "@"
$wsAmsin.Range("A54").Value = "2022-09-20"
$wsAmsin.Range("B54").Value = 44824.36937221065
$wsAmsin.Range("C54").NumberFormat = "@"
$wsAmsin.Range("C54").Value = "finalrun167"
$wsAmsin.Range("D54").Value = 75
$wsAmsin.Range("E54").Value = 72
$wsAmsin.Range("F54").Value = 3
$wsAmsin.Range("G54").Value = 1.89

$wsAmsin.Range("A51:G51").Copy()
$wsAmsin.Range("A52:G54").PasteSpecial(-4122)

# --- Sheet "BETA": one new row (28) ------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Range("A28").NumberFormat = "@"
$wsBeta.Range("A28").Value = "2022-09-20"
$wsBeta.Range("B28").Value = 44824.51490568312
$wsBeta.Range("C28").NumberFormat = "@"
$wsBeta.Range("C28").Value = "beta167"
$wsBeta.Range("D28").Value = 75
$wsBeta.Range("E28").Value = 75
$wsBeta.Range("F28").Value = 0
$wsBeta.Range("G28").Value = 2.12

$wsBeta.Range("A27:G27").Copy()
$wsBeta.Range("A28:G28").PasteSpecial(-4122)

# --- Sheet "AMS": restyle row 30, add new row 31 -----------------------
$wsAms = $wb.Worksheets.Item("AMS")

# row 30 existed but had lost its row formatting; rewrite it (text columns
# pinned to "@" first so the date-shaped text isn't reinterpreted as a
# date serial) with the Run Time value recalculated a hair more precisely,
# then restore the row's look from row 29 before appending the new
# Beta-server run underneath it.
$wsAms.Range("A30").NumberFormat = "@"
$wsAms.Range("A30").Value = "2022-08-24"
$wsAms.Range("B30").Value = 44797.91305877315
$wsAms.Range("C30").NumberFormat = "@"
$wsAms.Range("C30").Value = "166_live"
$wsAms.Range("D30").Value = 75
$wsAms.Range("E30").Value = 75
$wsAms.Range("F30").Value = 0
$wsAms.Range("G30").Value = 2.19

$wsAms.Range("A29:G29").Copy()
$wsAms.Range("A30:G30").PasteSpecial(-4122)

$wsAms.Range("A31").NumberFormat = "@"
$wsAms.Range("A31").Value = "2022-09-14"
$wsAms.Range("B31").Value = 44818.4146512037
$wsAms.Range("C31").NumberFormat = "@"
$wsAms.Range("C31").Value = "livehtfxsep166"
$wsAms.Range("D31").Value = 75
$wsAms.Range("E31").Value = 70
$wsAms.Range("F31").Value = 5
$wsAms.Range("G31").Value = 2.92

$wsAms.Range("A30:G30").Copy()
$wsAms.Range("A31:G31").PasteSpecial(-4122)

# Leave clipboard/selection state clean without disturbing which sheet
#/cell was active before this script ran.
$excel.CutCopyMode = 0
